# Applies the "Generate Report for Handoff" update:
#  - Overview sheet: status text "Handed back: in sync with en-US" -> "Ready for handoff"
#    and the "Latest HO Xliff Generate Date" timestamp refreshed.
#  - zh-cn / de-de sheets: Status -> "Ready for handoff", Priority "ht" -> "mt",
#    the handoff timestamp refreshed, and a new "handback file is not the
#    latest" error message recorded for the 66936f4b-... row.
#  - Column width tweaks on the three sheets to better fit the report.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8edf9892d6af96aa4f796d572767bed47a1ced36/e2e/66936f4b-add1-4a72-89ed-a1f62cca8d11.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bdda02b356990b467387a3f8f224780697132161/e2e/66936f4b-add1-4a72-89ed-a1f62cca8d11.md."

# ---- Overview sheet -----------------------------------------------------
$ws_overview.Range("E2:F3").Value = "Ready for handoff"
$ws_overview.Range("G2:G3").Value = "2016-11-14 07:29:34"

$ws_overview.Range("E1:F1").ColumnWidth = 16.333333333333332

# ---- zh-cn sheet ----------------------------------------------------------
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_zhcn.Range("C3").Value = "Ready for handoff"
$ws_zhcn.Range("E2:E3").Value = "mt"
$ws_zhcn.Range("H2:H3").Value = "2016-11-14 07:29:21"
$ws_zhcn.Range("P2").Value = $errorMessage

$ws_zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws_zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet ----------------------------------------------------------
$ws_dede.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("C3").Value = "Ready for handoff"
$ws_dede.Range("E2:E3").Value = "mt"
$ws_dede.Range("H2:H3").Value = "2016-11-14 07:29:34"
$ws_dede.Range("P2").Value = $errorMessage

$ws_dede.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws_dede.Columns.Item(16).ColumnWidth = 39.166666666666664
